$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# recall (row 3)
$ws.Range("F3").Value = 1

# f1-score (row 4)
$ws.Range("F4").Value = 0.2857142857142857

# f2-score (row 5)
$ws.Range("F5").Value = 0.5

# NDCG (row 6)
$ws.Range("F6").Value = 0.52129602861432

# M3 (row 8) - boolean
$ws.Range("F8").Value = $true

# M5 (row 9) - boolean
$ws.Range("F9").Value = $true

# position (row 10) - was empty inline string, now a number
$ws.Range("F10").Value = 2

# length (x of gs) (row 11)
$ws.Range("C11").Value = 5
$ws.Range("F11").Value = 6
$ws.Range("I11").Value = 2
$ws.Range("M11").Value = 4
$ws.Range("Q11").Value = 6
$ws.Range("S11").Value = 3
$ws.Range("U11").Value = 1
$ws.Range("W11").Value = 3
